$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Change 1: B5 was stored as text "1000274330" -> convert to a real number ---
$ws.Cells.Item(5, 2).Value = 1000274330

# --- Change 2: append new log row 6 ---
$ws.Cells.Item(6, 1).Value = "2025-10-15 22:24:54"

# B6 must stay TEXT ("1000135120"), not auto-coerced to a number.
# Leading apostrophe forces text entry; resetting the style afterwards
# strips the transient quote-prefix formatting Excel applies.
$ws.Cells.Item(6, 2).Value = "'1000135120"
$ws.Cells.Item(6, 2).Style = "Normal"

$ws.Cells.Item(6, 3).Value = "Leidy"
$ws.Cells.Item(6, 4).Value = "TARJETA DE CRÉDITO"

# E6 stays blank but present: touching a no-op formatting property keeps
# Excel from dropping the empty cell entirely.
$ws.Cells.Item(6, 5).Font.Bold = $false

$ws.Cells.Item(6, 6).Value = "PRORROGA CON PAGO"
$ws.Cells.Item(6, 7).Value = "48 cuotas"
$ws.Cells.Item(6, 8).Value = "34.19.100.134"
$ws.Cells.Item(6, 9).Value = "The Dalles"
$ws.Cells.Item(6, 10).Value = "Oregon"
$ws.Cells.Item(6, 11).Value = "United States"

# L6:P6 stay blank but present, same trick as E6.
$ws.Cells.Item(6, 12).Font.Bold = $false
$ws.Cells.Item(6, 13).Font.Bold = $false
$ws.Cells.Item(6, 14).Font.Bold = $false
$ws.Cells.Item(6, 15).Font.Bold = $false
$ws.Cells.Item(6, 16).Font.Bold = $false
